# Estadisticos Matutinos 15 Oct
# Updates the 2nd-partial and final statistics rows, and refreshes the
# "Rescatables" (remedial) roster with the current list of students.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Estadisticos 1P" -- update Aprobados/Por_Apro/Promedio columns
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 14
$ws1.Range("F2").Value = 29
$ws1.Range("G2").Value = 67.44
$ws1.Range("H2").Value = 7.4

$ws1.Range("D3").Value = 9
$ws1.Range("F3").Value = 20
$ws1.Range("G3").Value = 68.97
$ws1.Range("H3").Value = 7.7

# ---------------------------------------------------------------------
# Sheet "Estadisticos 2P" -- update Blancos column
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("E2").Value = 29
$ws2.Range("E3").Value = 20

# ---------------------------------------------------------------------
# Sheet "Estadisticos Final" -- update Aprobados/Por_Apro/Promedio columns
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 14
$ws3.Range("F2").Value = 29
$ws3.Range("G2").Value = 67.44
$ws3.Range("H2").Value = 7.4

$ws3.Range("D3").Value = 9
$ws3.Range("F3").Value = 20
$ws3.Range("G3").Value = 68.97
$ws3.Range("H3").Value = 7.7

# ---------------------------------------------------------------------
# Sheet "Rescatables" -- replace the roster with 5 students
# (write column-by-column so new shared strings land in the same order
# as the reference workbook: all Paterno, then all Materno, then Nombres)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$paterno = @("MIXCOHUA", "RAMIREZ", "DE LOS SANTOS", "PELLICO", "RUGERIO")
$materno = @("IXMATLAHUA", "PEREZ", "QUIÑONES", "TEQUIHUATLE", "SANCHEZ")
$nombres = @("ANGELINA", "LUIS REY", "GABRIELA", "JAZMIN", "KIMBERLY")
$nc      = @(21330051920088, 21330051920099, 21330051920104, 21330051920097, 21330051920383)

for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $ws4.Cells.Item($row, 2).Value = $paterno[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $ws4.Cells.Item($row, 3).Value = $materno[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $ws4.Cells.Item($row, 4).Value = $nombres[$i]
}
for ($i = 0; $i -lt 5; $i++) {
    $row = 2 + $i
    $ws4.Cells.Item($row, 1).Value = $nc[$i]
    $ws4.Cells.Item($row, 5).Value = "ÁLGEBRA"
    $ws4.Cells.Item($row, 6).Value = "1CV"
    $ws4.Cells.Item($row, 7).Value = 6
}
